# Update the "OKLGAtomicSwapInstance" row's USD price from 2 to 5.
# All the other cells in that row (D9:K9) are formulas referencing C9,
# so they will recalculate automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C9").Value = 5

$excel.CalculateFullRebuild()
